# Auto-generated edit script: updates crypto price/volume table
# to match the refreshed data from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "1.007").
# Force the whole Price column to Text format first so Excel keeps
# these as literal strings instead of silently parsing them into
# floating point numbers (which would corrupt values like "1.007").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.348.01"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.564.05"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").Value = "1.006"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "288.57"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "0.3732"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "49.21"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "0.3360"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "0.07413"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("D11").Value = "1.113"
$ws.Range("E11").Value = "  -4.29%  "
$ws.Range("D12").Value = "1.008"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "20.67"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").Value = "5.853"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "6.843"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "1.562.26"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "0.00001103"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "88.80"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "0.06696"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "6.123"
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").Value = "16.18"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").Value = "11.77"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").Value = "22.345.56"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "2.367"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").Value = "2.518"
$ws.Range("E26").Value = "  -10.62%  "
$ws.Range("D27").Value = "19.81"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("D28").Value = "146.54"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").Value = "4.993"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "124.16"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").Value = "1.733.93"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "1.990"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "0.9801"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D34").Value = "5.863"
$ws.Range("D35").Value = "9.619"
$ws.Range("E35").Value = "  -3.78%  "
$ws.Range("D36").Value = "0.08421"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "1.376"
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("D38").Value = "0.02434"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("D39").Value = "0.2250"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("D40").Value = "0.06355"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "5.305"
$ws.Range("E41").Value = "  -3.67%  "
$ws.Range("D42").Value = "0.6135"
$ws.Range("E42").Value = "  -3.07%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "1.005"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "10.85"
$ws.Range("E44").Value = "  -6.87%  "
$ws.Range("D45").Value = "13.78"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("D46").Value = "3.772"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").Value = "0.5729"
$ws.Range("E47").Value = "  -3.91%  "
$ws.Range("D48").Value = "2.021"
$ws.Range("E48").Value = "  -3.10%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "125.22"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "1.225"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("D51").Value = "0.07288"
$ws.Range("E51").Value = "  +0.43%  "

# Restore the default (unstyled) cell style on column D now that the
# values are locked in as text, so no stray number-format style is
# left behind on these cells.
$ws.Range("D2:D51").Style = "Normal"

